# "correções finais + quero morrer"
# Final corrections to capacity values (column E) on Tabelle1, plus the
# view-state left over from editing (frozen-pane scroll position / last
# selected cell) that Excel persists into the sheetView on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Capacity ("capacity", column E) corrections -------------------------
# Rows 12-19: capacity 5 -> 4
$ws.Range("E12:E19").Value = 4

# Row 40: capacity 5 -> 3
$ws.Range("E40").Value = 3

# Row 42: capacity 7 -> 3
$ws.Range("E42").Value = 3

# --- View state ------------------------------------------------------------
# Leave the workbook scrolled/selected where the author ended up: frozen
# top row, scrolled so row 2 is the first visible row below the freeze,
# with E38 as the active cell in the lower (frozen) pane.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("E38").Select() | Out-Null
